# Rename arc to link
#
# The "arcs" worksheet (2nd sheet, sheetId=1, r:id=rId2) is renamed to
# "links". It also becomes the active sheet/tab (activeTab moves from the
# "nodes" sheet to this one), with the selection on this sheet set to C36.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("arcs")
$ws.Name = "links"

# Make the renamed sheet the active tab and restore its selection.
$ws.Activate()
$ws.Range("C36").Select()
